$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reverse the "Periodo Mora" values in E16:E27 (previous EC periods replaced
# with the new set of periods, in the new/reverse order).
$periodos = @("2311","2310","2309","2308","2307","2306","2305","2304","2303","2302","2301","2212")
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}

# "Valor Mora" (F column) follows the period it belongs to: period 2212 pays
# 90000, period 2311 pays 84000 - these two values swap rows along with the
# reordering above.
$ws.Range("F16").Value = 84000
$ws.Range("F27").Value = 90000
